$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "af2d19bdad5bcc3708bcb2b9f32ccaf49E6BF287EE7E19ED8F9ADA2F466A0CBD4C488B69"
$ws.Range("H4").Value = "41454459d97f26fb5c2f8815b477a754FBAA0E440484A2E6480F313AAE783C9595E20036"
